$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 1).Value = "1000-1002"
$ws.Cells.Item(2, 2).Value = "9x5m"
$ws.Cells.Item(2, 3).Value = "Leif Wallén"
$ws.Cells.Item(2, 4).Value = "Plats"
$ws.Cells.Item(2, 5).Value = "15-04-2013 13:48:05"
# F2 stays an empty placeholder cell (no value), matching the original record.
$ws.Cells.Item(2, 6).Font.Bold = $true

$ws.Cells.Item(3, 1).Value = "1713-1715"
$ws.Cells.Item(3, 2).Value = "9x5 m"
$ws.Cells.Item(3, 3).Value = "Grilltösen"
$ws.Cells.Item(3, 4).Value = "Gatuköksprodukter: olika sorters korv, hamburgare, pommes frites"
$ws.Cells.Item(3, 5).Value = "02-07-2013 08:16:28"
$ws.Cells.Item(3, 6).Value = "asdasdasd"

$ws.Cells.Item(4, 1).Value = "1227-12228"
$ws.Cells.Item(4, 2).Value = "6x5 m"
$ws.Cells.Item(4, 3).Value = "Nightmare on tour AB"
$ws.Cells.Item(4, 4).Value = "asdasd"
$ws.Cells.Item(4, 5).Value = "02-07-2013 08:16:37"
$ws.Cells.Item(4, 6).Value = "asdasd"

$ws.Cells.Item(5, 1).Value = "1113-1115"
$ws.Cells.Item(5, 2).Value = "9x5 m"
$ws.Cells.Item(5, 3).Value = "Marknadsmedia"
$ws.Cells.Item(5, 4).Value = "asdasd"
$ws.Cells.Item(5, 5).Value = "02-07-2013 08:16:47"
$ws.Cells.Item(5, 6).Value = "asdasd"
